$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 26 (pushes the "ID 4 Ends" block, and
#    everything after it, down by one row). Excel inherits the row-25
#    formatting on insert, so normalise the new row's styles afterwards.
# ---------------------------------------------------------------------------
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = 42806
$ws.Range("A26").Style = "Normal"
$ws.Range("A26").NumberFormat = "mm/dd/yy;@"

$ws.Range("B26").Value = """Showcase"" build completed (focus on UI, user experience, aesthetics)"
$ws.Range("B26").Style = "Normal"

$ws.Range("C26:E26").Style = "Normal"
$ws.Range("C26:E26").ClearContents()

# ---------------------------------------------------------------------------
# 2. Fill in new Comments / Completed-by-end-of-ID columns for the rows in
#    the ID3 block (rows 14-24), and fix up the one milestone-name wording
#    change.
# ---------------------------------------------------------------------------

# Row 14 - Image-viewing room has a ceiling, no skybox visible
$ws.Range("D14").Value = "NO"
$ws.Range("D14").Style = "Bad"
$ws.Range("E14").Value = "Will be implemented when VR distance has been more formally calculated, to prevent objects clipping into the ceiling."

# Row 15 - File system interface implemented
$ws.Range("D15").Value = "NO"
$ws.Range("D15").Style = "Bad"
$ws.Range("E15").Value = "Exists in prototype form only"

# Row 16 - DICOM Library integrated into system
$ws.Range("E16").Value = "Fellow Oak DICOM library was selected for the project and crude tests pass"

# Row 17 - Risk Scan for ID3 completed
$ws.Range("D17").Value = "YES"
$ws.Range("D17").Style = "Good"

# Row 18 - Anti-motion-sickness report completed
$ws.Range("D18").Value = "YES"
$ws.Range("D18").Style = "Good"

# Row 19 - 40+ defects logged
$ws.Range("D19").Value = "YES"
$ws.Range("D19").Style = "Good"

# Row 21 - 3+ system tests implemented into project
$ws.Range("D21").Value = "YES"
$ws.Range("D21").Style = "Good"
$ws.Range("E21").Value = "Tests are implemented but still need development work to pass."

# Row 22 - Quit button wording tweak + Completed column
$ws.Range("B22").Value = "Quit button, other user-facing UI implemented"
$ws.Range("D22").Value = "YES"
$ws.Range("D22").Style = "Good"

# Row 23 - At least one image modification available for Copy objects
$ws.Range("D23").Value = "YES"
$ws.Range("D23").Style = "Good"
$ws.Range("E23").Value = "The Brightness slider was implemented, but may need some rework in a later sprint."

# Row 24 - 2+ user testing sessions completed (eg. users with glasses and no glasses)
$ws.Range("C24").Value = "NO"
$ws.Range("C24").Style = "Bad"
$ws.Range("D24").Value = "NO"
$ws.Range("D24").Style = "Bad"
$ws.Range("E24").Value = "Developers and testers agreed the build needed more work before formal user testing would be productive."

# ---------------------------------------------------------------------------
# 3. Re-sequence the milestones between the new row and the "ID 4 Ends"
#    divider (row 33 post-insert): two items drop off the list (moved /
#    superseded) and two new items take their place, so dates and labels
#    for rows 27-32 all need to be re-pointed at their new slot.
# ---------------------------------------------------------------------------

$ws.Range("A27").Value = 42810
$ws.Range("B27").Value = "Snap-to-grid for copies implemented"

$ws.Range("A28").Value = 42811
$ws.Range("B28").Value = "3+ medical professionals have demoed the software"

$ws.Range("A29").Value = 42811
$ws.Range("B29").Value = "1+ bug party completed"

$ws.Range("A30").Value = 42812
$ws.Range("B30").Value = "5+ user testing sessions completed"

$ws.Range("A31").Value = 42812
$ws.Range("B31").Value = "75+ defects logged"

$ws.Range("A32").Value = 42813
$ws.Range("B32").Value = "Coverage testing scheme correctly implemented"

# ---------------------------------------------------------------------------
# 4. Restore the final selection / active cell as recorded in the workbook.
# ---------------------------------------------------------------------------
$ws.Range("G20").Select()
